# Data source corrected and updated
#
# The source workbook recomputed columns J and K (rows 1-51) with new
# values: J -> 0.3 and K -> 0.5 for every row (row 1 previously held the
# text labels "r"/"s" via shared strings; those are replaced by the same
# numeric values used elsewhere in the column, so the shared-strings
# table becomes empty).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the full J and K columns (rows 1-51) to their corrected values.
$ws.Range("J1:J51").Value = 0.3
$ws.Range("K1:K51").Value = 0.5

# Reflect the author's final selection/view state: column K selected,
# anchored at K1.
$null = $ws.Range("K1:K51").Select()
